$wb = $excel.ActiveWorkbook

# --- Rename the first worksheet ---
# ("io_md_codes_phase_202304062020_" -> "io_md_codes_phase")
# This also causes the dependent defined name ("ExternalData_1"), which
# refers to the sheet by name, to be updated automatically.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "io_md_codes_phase"

# --- Move the active selection on that worksheet ---
# (was E46, scrolled so that A31 was the top-left cell; now K33)
$ws1.Activate()
$ws1.Range("K33").Select()
